$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 10:22"

# --- Suiza (row 13): refreshed case numbers, no re-sort needed ---
$ws.Range("B13").Value = 21176
$ws.Range("C13").Value = 76
$ws.Range("D13").Value = 7298
$ws.Range("E13").Value = 13154
$ws.Range("F13").Value = 391
$ws.Range("G13").Value = 9
$ws.Range("H13").Value = 724

# --- Rows 27-29: Dinamarca's totals overtake Chequia & Chile, so the
# three countries re-sort (descending by Casos totales). Dinamarca moves
# up to row 27 with refreshed numbers; Chequia and Chile keep their own
# data but shift down one row each. ---
$ws.Range("A27").Value = "Dinamarca"
$ws.Range("B27").Value = 4647
$ws.Range("C27").Value = 278
$ws.Range("D27").Value = 1327
$ws.Range("E27").Value = 3141
$ws.Range("F27").Value = 142
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 179

$ws.Range("A28").Value = "Chequia"
$ws.Range("B28").Value = 4591
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 96
$ws.Range("E28").Value = 4423
$ws.Range("F28").Value = 84
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 72

$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 4471
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 618
$ws.Range("E29").Value = 3819
$ws.Range("F29").Value = 307
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 34

# --- Row 31: Polonia refreshed numbers (stays in place) ---
$ws.Range("B31").Value = 4201
$ws.Range("C31").Value = 99
$ws.Range("D31").Value = 162
$ws.Range("E31").Value = 3941
$ws.Range("F31").Value = 50
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = 98

# --- Row 67: Lituania refreshed active/recovered counts ---
$ws.Range("D67").Value = 8
$ws.Range("E67").Value = 821

# --- Rows 176-178: Laos's totals overtake Granada & Sudan, same
# re-sort pattern as rows 27-29 above. ---
$ws.Range("A176").Value = "Laos"
$ws.Range("B176").Value = 12
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 12
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

$ws.Range("A177").Value = "Granada"
$ws.Range("B177").Value = 12
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 0
$ws.Range("E177").Value = 12
$ws.Range("F177").Value = 2
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

$ws.Range("A178").Value = "Sudan"
$ws.Range("B178").Value = 12
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 2
$ws.Range("E178").Value = 8
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 2
